# Regenerate save_data to use K (column G) instead of Strike#.
# Update the "K" column (G) values for the affected rows with the
# newly-calculated counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 3
    4  = 2
    5  = 1
    6  = 0
    7  = 3
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    14 = 0
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
